# Renamed few transcripts. Updated the DataSheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$rbdRows = @(2,5,7,8,10,11,13,15,17,18,20,22,25,26,27,35,36,40,42,43,45,49,50,51,53,55,56,58,59,62,73,80,83,84)
foreach ($r in $rbdRows) {
    $ws.Range("D$r").Value = "T"
}

$studentRows = @(63,69,70)
foreach ($r in $studentRows) {
    $ws.Range("D$r").Value = "S"
}
